$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17: Swoop Bite
$ws.Range("A17").Formula = "=ROW()-2"
$ws.Range("B17").Value = "Swoop Bite"
$ws.Range("C17").Value = "Weapons"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = "None"
$ws.Range("F17").Value = "Dives towards the player for a bite, then retreats upward. Main attack"

# Update existing rows 15/16 descriptions (column F)
$ws.Range("F15").Value = "The zombie swings its decaying arms in a wide arc. Main attack"
$ws.Range("F16").Value = "The dog lunges and bites viciously at the player's legs and arms. Main attack"

# New row 18: Gnaw
$ws.Range("A18").Formula = "=ROW()-2"
$ws.Range("B18").Value = "Gnaw"
$ws.Range("C18").Value = "Weapons"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "None"
$ws.Range("F18").Value = "Gnaws at the player and increases attack speed by 5%."

# New row 19: Sweeping Claw
$ws.Range("A19").Formula = "=ROW()-2"
$ws.Range("B19").Value = "Sweeping Claw"
$ws.Range("C19").Value = "Weapons"
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = "None"
$ws.Range("F19").Value = "Swipes both front limbs across the ground. Main attack"

$ws.Range("I19").Select()
